# Module 19 - Taxation Current Tax: add a second copy of WSE19.4
# ("WSE19.4 (2)") with updated b/f figures, and rework the journal
# workings on the original WSE19.4 sheet (swap the finance-cost /
# current-tax-expense journal rows, add a new b/f-balance
# reconciliation block at the bottom of the sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# STEP 1: duplicate WSE19.4 -> "WSE19.4 (2)" (placed immediately
# after the original, exactly like Excel's own "Move or Copy..."
# with "Create a copy" ticked).
# ---------------------------------------------------------------
$wsOrig = $wb.Worksheets.Item("WSE19.4")
$wsOrig.Copy([System.Reflection.Missing]::Value, $wsOrig)
$wsCopy = $wb.Worksheets.Item("WSE19.4 (2)")

# ---------------------------------------------------------------
# STEP 2: on the NEW copy, bump the b/f creditor figures up to the
# restated amounts and leave the rest of the workings to recalc.
# ---------------------------------------------------------------
$wsCopy.Range("C34").Value = 1051600
$wsCopy.Range("C35").Value = 2103200
$wsCopy.Range("C43").Select()

# ---------------------------------------------------------------
# STEP 3: rework the original WSE19.4 sheet.
# ---------------------------------------------------------------
$wsOrig.Activate()

# Widen column C slightly so the new "cr - finance income" /
# "cr - current tax expense" labels aren't clipped.
$wsOrig.Columns.Item(3).ColumnWidth = 18.57

# New b/f-balance reconciliation block at the foot of the sheet.
# (entered before the rows-55/56 relabel below so new shared-string
# entries land in the same order the original author created them)
$wsOrig.Range("B84").Value = "b/f balance"
$wsOrig.Range("D84").Value = 177
$wsOrig.Range("B85").Value = "PY"
$wsOrig.Range("D85").Value = 170
$wsOrig.Range("D86").Formula = "=+D84-D85"
$wsOrig.Range("B88").Value = "the CY payments go through PL"
$wsOrig.Range("B91").Value = "overpaid tax - the expected too much in PY"

# Row 56 used to hold "cr - finance income" / +D74 in columns C/F;
# that now moves down to row 57, and row 56 becomes the new
# "dr - current tax expense" leg referencing the (relocated) F57.
$wsOrig.Range("C56").Clear()
$wsOrig.Range("F56").Clear()
$wsOrig.Range("B56").Value = "dr - current tax expense"
$wsOrig.Range("E56").Formula = "=+F57"

# Row 55: re-label from "total corporation tax receivable" to the
# first leg of the new journal; the D40 reference stays as-is.
$wsOrig.Range("B55").Value = "dr - finance cost"

# Row 57: now holds "cr - finance income" with its own +D74 total.
$wsOrig.Range("C57").Value = "cr - finance income"
$wsOrig.Range("F57").Formula = "=+D74"

# Row 58: keeps its "cr - current tax expense" label, but the
# total is now a plain entered amount rather than a -D76 formula.
$wsOrig.Range("C58").Value = "cr - current tax expense"
$wsOrig.Range("F58").Value = 80

# Downstream workings: point at the relocated F57, and rebuild the
# interest-on-late-payment total from its components instead of F58.
$wsOrig.Range("E63").Formula = "=+F57"
$wsOrig.Range("E66").Formula = "=+E62+E63-E64"

$wsOrig.Range("B52").Select()

# ---------------------------------------------------------------
# STEP 4: touch up the stored selections on the neighbouring
# worksheets.
# ---------------------------------------------------------------
$wsE2 = $wb.Worksheets.Item("WSE19.2")
$wsE2.Activate()
$wsE2.Range("D31").Select()

$wsE3 = $wb.Worksheets.Item("WSE19.3")
$wsE3.Activate()
$wsE3.Range("D39").Select()

# ---------------------------------------------------------------
# STEP 5: leave WSE19.4 as the active/selected tab, matching the
# original workbook's last-saved state.
# ---------------------------------------------------------------
$wsOrig.Activate()
